# Clear the quantile values (A:D) for the "Spices" row (row 10) and
# flip its reference-diet "type" from B to A -- mirrors adding the
# second reference diet and re-running the results for that row.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A10:D10").ClearContents()
$ws.Range("F10").Value = "A"
